$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.2035556666666667
$ws.Range("H2").Value2 = 0.6106670000000001
$ws.Range("I2").Value2 = 0.006148914270823412
$ws.Range("J2").Value2 = 0.006148914270823412
$ws.Range("M2").Value2 = 61.04160633333334
$ws.Range("N2").Value2 = 183.124819
$ws.Range("O2").Value2 = 0.2043613460574534
$ws.Range("P2").Value2 = 0.2043613460574534
$ws.Range("Q2").Value2 = 12.42536487158589
$ws.Range("R2").Value2 = 111.828283844273
$ws.Range("S2").Value2 = 0.001256600397177357
$ws.Range("T2").Value2 = 0.001256600397177357
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.2035556666666667
$ws.Range("H3").Value2 = 0.6106670000000001
$ws.Range("I3").Value2 = 0.006148914270823412
$ws.Range("J3").Value2 = 0.006148914270823412
$ws.Range("O3").Value2 = 0.3559304658284363
$ws.Range("P3").Value2 = 0.3559304658284363
$ws.Range("Q3").Value2 = 21.64091200294067
$ws.Range("R3").Value2 = 194.768208026466
$ws.Range("S3").Value2 = 0.002188585920753296
$ws.Range("T3").Value2 = 0.002188585920753297
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.2035556666666667
$ws.Range("H4").Value2 = 0.6106670000000001
$ws.Range("I4").Value2 = 0.006148914270823412
$ws.Range("J4").Value2 = 0.006148914270823412
$ws.Range("M4").Value2 = 131.3384093333333
$ws.Range("N4").Value2 = 394.015228
$ws.Range("O4").Value2 = 0.4397081881141102
$ws.Range("P4").Value2 = 0.4397081881141103
$ws.Range("Q4").Value2 = 26.73467747078622
$ws.Range("R4").Value2 = 240.612097237076
$ws.Range("S4").Value2 = 0.002703727952892758
$ws.Range("T4").Value2 = 0.002703727952892758
$ws.Range("I5").Value2 = 0.735846381812327
$ws.Range("J5").Value2 = 0.735846381812327
$ws.Range("M5").Value2 = 61.04160633333334
$ws.Range("N5").Value2 = 183.124819
$ws.Range("O5").Value2 = 0.2043613460574534
$ws.Range("P5").Value2 = 0.2043613460574534
$ws.Range("Q5").Value2 = 1486.955156756494
$ws.Range("R5").Value2 = 13382.59641080844
$ws.Range("S5").Value2 = 0.150378557078674
$ws.Range("T5").Value2 = 0.150378557078674
$ws.Range("I6").Value2 = 0.735846381812327
$ws.Range("J6").Value2 = 0.735846381812327
$ws.Range("O6").Value2 = 0.3559304658284363
$ws.Range("P6").Value2 = 0.3559304658284363
$ws.Range("S6").Value2 = 0.2619101454566309
$ws.Range("T6").Value2 = 0.261910145456631
$ws.Range("I7").Value2 = 0.735846381812327
$ws.Range("J7").Value2 = 0.735846381812327
$ws.Range("M7").Value2 = 131.3384093333333
$ws.Range("N7").Value2 = 394.015228
$ws.Range("O7").Value2 = 0.4397081881141102
$ws.Range("P7").Value2 = 0.4397081881141103
$ws.Range("Q7").Value2 = 3199.36411849879
$ws.Range("R7").Value2 = 28794.27706648911
$ws.Range("S7").Value2 = 0.3235576792770221
$ws.Range("T7").Value2 = 0.3235576792770221
$ws.Range("G8").Value2 = 8.541072
$ws.Range("H8").Value2 = 25.623216
$ws.Range("I8").Value2 = 0.2580047039168495
$ws.Range("J8").Value2 = 0.2580047039168495
$ws.Range("M8").Value2 = 61.04160633333334
$ws.Range("N8").Value2 = 183.124819
$ws.Range("O8").Value2 = 0.2043613460574534
$ws.Range("P8").Value2 = 0.2043613460574534
$ws.Range("Q8").Value2 = 521.360754688656
$ws.Range("R8").Value2 = 4692.246792197905
$ws.Range("S8").Value2 = 0.0527261885816021
$ws.Range("T8").Value2 = 0.0527261885816021
$ws.Range("G9").Value2 = 8.541072
$ws.Range("H9").Value2 = 25.623216
$ws.Range("I9").Value2 = 0.2580047039168495
$ws.Range("J9").Value2 = 0.2580047039168495
$ws.Range("O9").Value2 = 0.3559304658284363
$ws.Range("P9").Value2 = 0.3559304658284363
$ws.Range("Q9").Value2 = 908.039508747552
$ws.Range("R9").Value2 = 8172.355578727967
$ws.Range("S9").Value2 = 0.09183173445105204
$ws.Range("T9").Value2 = 0.09183173445105206
$ws.Range("G10").Value2 = 8.541072
$ws.Range("H10").Value2 = 25.623216
$ws.Range("I10").Value2 = 0.2580047039168495
$ws.Range("J10").Value2 = 0.2580047039168495
$ws.Range("M10").Value2 = 131.3384093333333
$ws.Range("N10").Value2 = 394.015228
$ws.Range("O10").Value2 = 0.4397081881141102
$ws.Range("P10").Value2 = 0.4397081881141103
$ws.Range("Q10").Value2 = 1121.770810481472
$ws.Range("R10").Value2 = 10095.93729433325
$ws.Range("S10").Value2 = 0.1134467808841954
$ws.Range("T10").Value2 = 0.1134467808841954
